$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the filter-description text that used to live in A1 (row 1 becomes empty / removed)
$ws.Range("A1").ClearContents()

# Update the turnover (ТО) value for the second salesperson (row 5)
$ws.Range("C5").Value = 7777

# Move the active selection, matching where the author's cursor ended up
$ws.Range("G16").Select()
